$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update ---
$ws.Range("A1").Value = 'Datos actualizados a 18 de Abril de 2020 a las 22:47'
$ws.Range("B4").Value = 734631
$ws.Range("C4").Value = 24896
$ws.Range("D4").Value = 67158
$ws.Range("E4").Value = 628700
$ws.Range("F4").Value = 13536
$ws.Range("G4").Value = 1619
$ws.Range("H4").Value = 38773
$ws.Range("B8").Value = 143172
$ws.Range("C8").Value = 1775
$ws.Range("E8").Value = 53320
$ws.Range("G8").Value = 100
$ws.Range("H8").Value = 4452
$ws.Range("B27").Value = 10296
$ws.Range("C27").Value = 509
$ws.Range("D27").Value = 1069
$ws.Range("E27").Value = 9005
$ws.Range("F27").Value = 221
$ws.Range("G27").Value = 32
$ws.Range("H27").Value = 222
$ws.Range("F53").Value = 36
$ws.Range("B94").Value = 655
$ws.Range("C94").Value = 6
$ws.Range("D94").Value = 97
$ws.Range("E94").Value = 554
$ws.Range("D149").Value = 17
$ws.Range("E149").Value = 53
$ws.Range("A167").Value = 'Maldivas'
$ws.Range("B167").Value = 35
$ws.Range("C167").Value = 6
$ws.Range("D167").Value = 16
$ws.Range("E167").Value = 19
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 0
$ws.Range("A168").Value = 'Benin'
$ws.Range("B168").Value = 35
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 18
$ws.Range("E168").Value = 16
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 1
$ws.Range("A169").Value = 'Mozambique'
$ws.Range("B169").Value = 34
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 2
$ws.Range("E169").Value = 32
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 0
